# atualizacao matriz de risco
# Renumber the "numero_ata" values (column AK) from the 163-174 series
# to the 055-066 series, keeping the same relative order/mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "785810/2024-163/00" = "785810/2024-055/00"
    "785810/2024-164/00" = "785810/2024-056/00"
    "785810/2024-165/00" = "785810/2024-057/00"
    "785810/2024-166/00" = "785810/2024-058/00"
    "785810/2024-167/00" = "785810/2024-059/00"
    "785810/2024-168/00" = "785810/2024-060/00"
    "785810/2024-169/00" = "785810/2024-061/00"
    "785810/2024-170/00" = "785810/2024-062/00"
    "785810/2024-171/00" = "785810/2024-063/00"
    "785810/2024-172/00" = "785810/2024-064/00"
    "785810/2024-173/00" = "785810/2024-065/00"
    "785810/2024-174/00" = "785810/2024-066/00"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 37).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 37)
    $val = $cell.Value2
    if ($val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
